$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the locked cells below
# (which hold the model's weights/percent-change figures and the
# disclosure footnote) can be updated, then restore protection afterward.
$ws.Unprotect()

# Footnote in A16: bump the "as of" model date from 2021-05-10 to 2021-05-11.
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."

# Refreshed Weight (D) / Percent Change (E) figures for each holding row.
$updates = @(
    @{Row=2;  D=0.03101670610048274;     E=-0.01256366723259761},
    @{Row=3;  D=0.02398263949618754;     E=-0.01205936920222628},
    @{Row=4;  D=0.05239747109205833;     E=-0.009675190048375804},
    @{Row=5;  D=0.137571318720929;       E=-0.007478458787189113},
    @{Row=6;  D=0.03206493479752429;     E=-0.02208419599723954},
    @{Row=7;  D=0.1191182098554325;      E=-0.01160736815543784},
    @{Row=8;  D=0.1038444598072102;      E=-0.01296362981634858},
    @{Row=9;  D=0.02962672029381645;     E=0.004277856997351615},
    @{Row=10; D=0.1273068900742044;      E=-0.01635991820040905},
    @{Row=11; D=0.2409837303261863;      E=-0.0006531069229332997},
    @{Row=12; D=0.1020869194359682;      E=-0.00292112950340806},
    @{Row=13; D=1;                       E=-0.008063243052496727}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Restore sheet protection (same basic options as before the edit).
$ws.Protect()
